# 1. Remove the obsolete "LE-ADM-Arbeitsbereich1" row (row 3) from Bureau and
#    LienBureauLieuGestion; the rows below shift up and the Bureau-id column
#    is renumbered by the workbook's own formulas/values already in place.
$wb = $excel.ActiveWorkbook

$wsBureau = $wb.Worksheets.Item("Bureau")
$wsBureau.Rows.Item(3).Delete()

$wsLien = $wb.Worksheets.Item("LienBureauLieuGestion")
$wsLien.Rows.Item(3).Delete()

# 2. Add the new "BureauUsers" sheet at the end of the workbook (TODO: used
#    later to default users to their Bureau's Windows printer).
$sheetCount = $wb.Worksheets.Count
$wsUsers = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$wsUsers.Name = "BureauUsers"

$wsUsers.Cells.Item(1, 1).Value = "root.Profiles.Bureau-id"
$wsUsers.Cells.Item(1, 2).Value = "root.Profiles.Bureau-libelle"
$wsUsers.Cells.Item(1, 3).Value = "users"

$bureauUsersData = @(
    ,@(1, 'LE-AAU-Arbeitsbereich 1', '{''B126CHK'', ''B126GRS'', ''B126GUM'', ''B126VOK'', ''B126PAF'', ''B126MA7'', ''B126CHR'', ''B126PAS'', ''B126AMA'', ''B126BC1'', ''B126JAG'', ''B126SIS'', ''B126ISM''}')
    ,@(1, 'AL-ZUL-CC61', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(2, 'AL-ZUL-PEZ 1', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(3, 'AL-ZUL-PEZ 2', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(4, 'AL-ZUL-PEZ 3', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(5, 'AL-ZUL-PEZ 4', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(6, 'AL-ZUL-PEZ 5', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(7, 'AL-ZUL-PEZ 6', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(8, 'AL-ZUL-PEZ 7', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(9, 'AL-ZUL-FZZ 1', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(10, 'AL-ZUL-FZZ 2', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(11, 'AL-ZUL-FZZ 3', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(12, 'AL-ZUL-FZZ 4', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(13, 'AL-ZUL-FZZ 5', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(14, 'AL-ZUL-FZZ 6', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(15, 'AL-ZUL-FZZ 7', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(16, 'AL-ZUL-FZZSpez1', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(17, 'AL-ZUL-FZZSpez2', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(18, 'AL-ZUL-FZZSpez3', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(19, 'AL-ZUL-FZZSpez4', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(20, 'AL-ZUL-FZZSpez5', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(21, 'AL-ZUL-FZZSpez6', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(22, 'AL-ZUL-FZZSpez7', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(23, 'AL-ZUL-FZZSpez8', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(24, 'AL-ZUL-FZZSpez9', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(25, 'AL-ZUL-FZZSpez10', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(26, 'AL-ZUL-FZZSpez11', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(27, 'AL-ZUL-FZZSpez12', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(28, 'AL-ZUL-FZZSpez13', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(29, 'AL-ZUL-FZZSpez14', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(30, 'AL-ZUL-FZZSpez15', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(31, 'AL-ZUL-FZZSpez16', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(32, 'AL-ZUL-FZZSpez17', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(33, 'AL-ZUL-FZZSpez18', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(34, 'AL-ZUL-FZZSpez19', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
    ,@(35, 'AL-ZUL-FZZSpez20', '{''B126FRC'', ''B126GRG'', ''B126MAS'', ''B126IMD'', ''B126SMP'', ''B126LOG'', ''B126HNA''}')
)

$r = 2
foreach ($row in $bureauUsersData) {
    $wsUsers.Cells.Item($r, 1).Value = $row[0]
    $wsUsers.Cells.Item($r, 2).Value = $row[1]
    $wsUsers.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
